# Reset every contributor's commit count to 0, keeping the "Name : N" text format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2"  = "LeeYinWei : 0"
    "C2"  = "unknown899 : 0"
    "B3"  = "neoAurora : 0"
    "D3"  = "howardhung14 : 0"
    "B4"  = "yoyo0213 : 0"
    "C4"  = "JonathanYangSW : 0"
    "D4"  = "GinoChen113511247 : 0"
    "C5"  = "peienwu1216 : 0"
    "D5"  = "chxyuuu : 0"
    "B6"  = "ginny923 : 0"
    "C6"  = "joanna0420 : 0"
    "D6"  = "dua0505 : 0"
    "B7"  = "jui-pixel : 0"
    "D7"  = "charles691 : 0"
    "B8"  = "Tony104147 : 0"
    "B9"  = "haleychang0530 : 0"
    "C9"  = "Hazel-1212 : 0"
    "B10" = "CHENG-JE : 0"
    "C10" = "lwc-ed : 0"
    "B11" = "tpvupu : 0"
    "C11" = "xiaotin22 : 0"
    "D11" = "calistayang : 0"
    "B13" = "kufanghua : 0"
    "C13" = "yezh0915 : 0"
    "D13" = "fiesta0217 : 0"
    "D14" = "jing1688 : 0"
    "B15" = "weiouo-0817 : 0"
    "B16" = "gamemode0701 : 0"
    "C16" = "Tonyyu2403 : 0"
    "B17" = "TerryCheese : 0"
    "C17" = "junlin27 : 0"
    "B18" = "Miiaow3011 : 0"
    "C18" = "bonnieliao774 : 0"
    "D18" = "emmazheng0318 : 0"
    "C19" = "TedChueh : 0"
    "B20" = "max052028 : 0"
    "B21" = "houyuankai : 0"
    "B22" = "0u88 : 0"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
